{"js": "// Replace each two-digit multiplication equation in the document with its\n// updated value. Every equation text is unique within the document, so we\n// can safely search for the old text and replace it with the new text.\nconst replacements = [\n  [\"72\u00d770=\", \"41\u00d745=\"],\n  [\"57\u00d789=\", \"36\u00d740=\"],\n  [\"60\u00d756=\", \"34\u00d743=\"],\n  [\"45\u00d737=\", \"43\u00d740=\"],\n  [\"68\u00d761=\", \"90\u00d780=\"],\n  [\"67\u00d723=\", \"22\u00d764=\"],\n  [\"71\u00d721=\", \"61\u00d792=\"],\n  [\"33\u00d775=\", \"61\u00d726=\"],\n  [\"37\u00d732=\", \"55\u00d786=\"],\n  [\"76\u00d775=\", \"30\u00d717=\"],\n  [\"97\u00d781=\", \"95\u00d746=\"],\n  [\"67\u00d760=\", \"15\u00d738=\"],\n  [\"86\u00d715=\", \"48\u00d721=\"],\n  [\"90\u00d718=\", \"44\u00d794=\"],\n  [\"55\u00d773=\", \"91\u00d780=\"],\n  [\"71\u00d791=\", \"39\u00d748=\"],\n  [\"82\u00d750=\", \"48\u00d731=\"],\n  [\"23\u00d787=\", \"35\u00d731=\"],\n  [\"48\u00d736=\", \"97\u00d735=\"],\n  [\"71\u00d759=\", \"90\u00d766=\"],\n  [\"38\u00d719=\", \"34\u00d797=\"],\n  [\"53\u00d726=\", \"49\u00d745=\"],\n  [\"51\u00d749=\", \"92\u00d781=\"],\n  [\"93\u00d795=\", \"77\u00d778=\"],\n  [\"20\u00d751=\", \"14\u00d738=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Replace each two-digit multiplication equation in the document with its\n# updated value. Every equation text is unique within the document, so a\n# simple Find/Replace per pair is safe and unambiguous.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"72\u00d770=\", \"41\u00d745=\"),\n    @(\"57\u00d789=\", \"36\u00d740=\"),\n    @(\"60\u00d756=\", \"34\u00d743=\"),\n    @(\"45\u00d737=\", \"43\u00d740=\"),\n    @(\"68\u00d761=\", \"90\u00d780=\"),\n    @(\"67\u00d723=\", \"22\u00d764=\"),\n    @(\"71\u00d721=\", \"61\u00d792=\"),\n    @(\"33\u00d775=\", \"61\u00d726=\"),\n    @(\"37\u00d732=\", \"55\u00d786=\"),\n    @(\"76\u00d775=\", \"30\u00d717=\"),\n    @(\"97\u00d781=\", \"95\u00d746=\"),\n    @(\"67\u00d760=\", \"15\u00d738=\"),\n    @(\"86\u00d715=\", \"48\u00d721=\"),\n    @(\"90\u00d718=\", \"44\u00d794=\"),\n    @(\"55\u00d773=\", \"91\u00d780=\"),\n    @(\"71\u00d791=\", \"39\u00d748=\"),\n    @(\"82\u00d750=\", \"48\u00d731=\"),\n    @(\"23\u00d787=\", \"35\u00d731=\"),\n    @(\"48\u00d736=\", \"97\u00d735=\"),\n    @(\"71\u00d759=\", \"90\u00d766=\"),\n    @(\"38\u00d719=\", \"34\u00d797=\"),\n    @(\"53\u00d726=\", \"49\u00d745=\"),\n    @(\"51\u00d749=\", \"92\u00d781=\"),\n    @(\"93\u00d795=\", \"77\u00d778=\"),\n    @(\"20\u00d751=\", \"14\u00d738=\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    $find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n}\n"}
